$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to round-trip a cell's existing format (fill/font/alignment/
# number-format) across a text-forcing write, so Student ID values that look
# numeric ("211131") stay stored as text (matching the source file's
# t="inlineStr" cells) instead of Excel's default "looks like a number" ->
# numeric-cell inference, while the original cell style index is preserved.
$scratch = $ws.Cells.Item(1, 8)

function Set-TextValue($row, $value) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Copy()
    $scratch.PasteSpecial(-4122)  # xlPasteFormats
    $cell.Value = "'" + $value    # leading apostrophe forces text entry
    $scratch.Copy()
    $cell.PasteSpecial(-4122)     # restore the original cell format/style
    $scratch.Clear()
}

Set-TextValue 2  "211131"
Set-TextValue 3  "211043"
Set-TextValue 4  "191375"
Set-TextValue 5  "201574"
Set-TextValue 6  "211010"
Set-TextValue 7  "201080"
Set-TextValue 8  "201465"
Set-TextValue 9  "190801"
Set-TextValue 10 "200914"
Set-TextValue 11 "200938"
Set-TextValue 12 "200850"
Set-TextValue 13 "200877"
Set-TextValue 14 "202089"
Set-TextValue 15 "201834"
Set-TextValue 16 "200491"
Set-TextValue 17 "201840"
Set-TextValue 18 "190975"
Set-TextValue 19 "201825"
Set-TextValue 20 "201255"
Set-TextValue 21 "210728"
Set-TextValue 22 "201682"
Set-TextValue 23 "201397"
Set-TextValue 24 "200344"
Set-TextValue 25 "200804"

# Drop the last two log rows (26-27) entirely.
$ws.Rows("26:27").Delete()
